$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("crisis_src")
$ws.Activate()

# Reset H3:H47 and Q3:Q47 to 0 (prior-month totals reset for rerun)
$ws.Range("H3:H47").Value = 0
$ws.Range("Q3:Q47").Value = 0

# Clear the derived "Proportion of Cross Systems Populations Served" column (S) for data rows
$ws.Range("S3:S29").ClearContents()

# Update the active selection to reflect where the user left off
$ws.Range("K41").Select()
